$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '80.381.95'
$ws.Range('E2').Value = '  +4.56%  '

# Row 3
$ws.Range('D3').Value = '3.188.81'
$ws.Range('E3').Value = '  +1.34%  '

# Row 4
$ws.Range('E4').Value = '  +0.15%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '209.23'
$ws.Range('E5').Value = '  +3.80%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '627.19'
$ws.Range('E6').Value = '  +0.26%  '

# Row 7
$ws.Range('E7').Value = '  +25.68%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.03%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.588'
$ws.Range('E9').Value = '  +5.02%  '

# Row 10
$ws.Range('D10').Value = '3.186.04'
$ws.Range('E10').Value = '  +1.25%  '

# Row 11
$ws.Range('E11').Value = '  +20.41%  '

# Row 12
$ws.Range('E12').Value = '  +26.49%  '

# Row 13
$ws.Range('E13').Value = '  +1.22%  '

# Row 14
$ws.Range('D14').Value = '3.771.42'
$ws.Range('E14').Value = '  +1.34%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.29'
$ws.Range('E15').Value = '  +0.40%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '31.79'
$ws.Range('E16').Value = '  +6.28%  '

# Row 17
$ws.Range('D17').Value = '80.414.06'
$ws.Range('E17').Value = '  +4.71%  '

# Row 18
$ws.Range('D18').Value = '3.184.70'
$ws.Range('E18').Value = '  +1.35%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.24'
$ws.Range('E19').Value = '  +3.32%  '

# Row 20
$ws.Range('E20').Value = '  +8.81%  '

# Row 21
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '440.23'
$ws.Range('E21').Value = '  +8.88%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.18'
$ws.Range('E22').Value = '  -1.38%  '

# Row 23
$ws.Range('E23').Value = '  +12.31%  '

# Row 24
$ws.Range('E24').Value = '  +5.86%  '

# Row 25
$ws.Range('D25').Value = '3.357.01'
$ws.Range('E25').Value = '  +1.45%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '76.24'
$ws.Range('E26').Value = '  +2.99%  '

# Row 27
$ws.Range('E27').Value = '  +0.75%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.87'
$ws.Range('E28').Value = '  +3.74%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.07%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0000122'
$ws.Range('E30').Value = '  +7.60%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  +0.46%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.94'
$ws.Range('E32').Value = '  +4.43%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '559.25'
$ws.Range('E33').Value = '  +6.96%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.44'
$ws.Range('E34').Value = '  -1.42%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.151'
$ws.Range('E35').Value = '  +11.45%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.99'
$ws.Range('E36').Value = '  +1.38%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '23.10'
$ws.Range('E37').Value = '  +6.40%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.124'
$ws.Range('E38').Value = '  +19.55%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  +0.33%  '

# Row 40
$ws.Range('E40').Value = '  +5.27%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '20.76'
$ws.Range('E41').Value = '  +3.41%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '162.90'
$ws.Range('E42').Value = '  -0.33%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.65'
$ws.Range('E43').Value = '  +5.24%  '

# Row 44
$ws.Range('E44').Value = '  +0.02%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '189.49'
$ws.Range('E45').Value = '  -3.68%  '

# Row 46
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.72'
$ws.Range('E46').Value = '  +8.61%  '

# Row 47
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.81'
$ws.Range('E47').Value = '  +3.92%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.786'
$ws.Range('E48').Value = '  -3.04%  '

# Row 49
$ws.Range('E49').Value = '  -0.01%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '42.90'
$ws.Range('E50').Value = '  +1.74%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.25'
$ws.Range('E51').Value = '  +5.05%  '

# Reset number format back to General style so cell styling matches original (no style index)
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Style = 'Normal'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
